$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, reusing the same cell formatting as the
# other header cells (bold font, border, centered alignment) by copying
# an existing header cell's format rather than rebuilding it by hand
# (keeps the same shared style entry instead of minting a new one).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add new numeric value in H2 (plain, unstyled like the other data cells)
$ws.Range("H2").Value = 0
